$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G duplicates column F ("status" header / "Pass" values)
$ws.Range("G1").Value = "status"
$ws.Range("G2").Value = "Pass"
$ws.Range("G3").Value = "Pass"
$ws.Range("G4").Value = "Pass"

# G1 picks up the same highlighted fill as F1 (new style record, like F1's)
$ws.Range("G1").Interior.ColorIndex = $ws.Range("F1").Interior.ColorIndex()

# Column G should be sized like column F
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth()
